$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Days remaining" values for two trials following a data refresh.
# Row 6 = REJOICE (MK-5909-003): 8 -> 7
# Row 8 = REMASTER (CLOU): 28 -> 27
$ws.Range("B6").Value = 7
$ws.Range("B8").Value = 27
